{"js": "// The paragraph \"Systemet gemmer oplysningerne.\" should become a top-level\n// bullet (ilvl 0 instead of 1), and it gains the \"_GoBack\" bookmark at its\n// very start (this bookmark previously sat, orphaned, inside the\n// \"Vinduet forts\u00e6tter fra usecase 02.\" paragraph, splitting one sentence\n// into two runs). We also normalize that split back into a single run.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet targetPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t && t.trim() === \"Systemet gemmer oplysningerne.\") {\n    targetPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (targetPara) {\n  // Promote the bullet from second level (ilvl 1) to first level (ilvl 0).\n  targetPara.listItem.level = 0;\n}\n\n// Remove the old \"_GoBack\" bookmark wherever it currently lives \u2026\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// \u2026 merge the two runs that used to straddle it back into a single run \u2026\nconst hits = context.document.body.search(\"Vinduet forts\u00e6tter fra \", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\"Vinduet forts\u00e6tter fra \", \"Replace\");\n}\n\n// \u2026 and re-insert \"_GoBack\" at the start of the promoted bullet paragraph.\nif (targetPara) {\n  const startRange = targetPara.getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# The paragraph \"Systemet gemmer oplysningerne.\" should become a top-level\n# bullet (ilvl 0 instead of 1, i.e. Word's 1-based ListLevelNumber 1 instead\n# of 2), and it gains the \"_GoBack\" bookmark at its very start. That bookmark\n# used to sit, orphaned, inside the \"Vinduet forts\u00e6tter fra usecase 02.\"\n# paragraph, splitting one sentence into two runs (\"...forts\u00e6tter fr\" / \"a \").\n# We normalize that split back into a single run too.\n\n$d = $word.ActiveDocument\n\n# Find the \"Systemet gemmer oplysningerne.\" paragraph.\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -eq \"Systemet gemmer oplysningerne.\") {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -ne $null) {\n    # Promote the bullet from second level (ilvl 1) to first level (ilvl 0);\n    # Word's ListLevelNumber is 1-based, so that is level 1.\n    $targetPara.Range.ListFormat.ListLevelNumber = 1\n}\n\n# Merge the two runs that used to straddle the stray bookmark back into a\n# single run; replacing found text with itself rebuilds it as one run and\n# also drops any bookmark boundary that fell inside the found range.\n$find = $d.Content\n$find.Find.ClearFormatting()\n$find.Find.Text = \"Vinduet forts\u00e6tter fra \"\n$find.Find.Replacement.ClearFormatting()\n$find.Find.Replacement.Text = \"Vinduet forts\u00e6tter fra \"\n$find.Find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2)\n\n# Make sure no stray \"_GoBack\" bookmark is left anywhere else \u2026\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# \u2026 then re-insert it collapsed at the very start of the promoted paragraph.\nif ($targetPara -ne $null) {\n    $bmRange = $targetPara.Range.Duplicate\n    $bmRange.Collapse(1)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
